# "Season 14, matchdays prepares"
# Остаток без Оксанич Кирилл: удаляем строку 15 (целиком), остальные
# строки автоматически сдвигаются вверх, а соответствующая неиспользуемая
# строка в sharedStrings.xml будет убрана движком сама.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 15 (Оксанич Кирилл) - rows below shift up.
$ws.Rows(15).Delete() | Out-Null

# Re-apply the existing sort (by column B, participant name) over the
# now-smaller data range so the sheet's sortState/AutoSort metadata is
# refreshed to match the new extent (A2:AA18).
$ws.Sort.SortFields.Clear() | Out-Null
$ws.Sort.SortFields.Add($ws.Range("B2:B18")) | Out-Null
$ws.Sort.SetRange($ws.Range("A2:AA18")) | Out-Null
$ws.Sort.Header = [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlNo
$ws.Sort.Apply() | Out-Null

# Match the saved selection from the edit: the whole 15th row selected.
$ws.Range("A15:XFD15").Select() | Out-Null
